$wb = $excel.ActiveWorkbook

$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Xiaojiawa Coal Mine, China, M0384, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 9; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
